$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 20.54580371028513
$ws.Range("C2").Value = 4.119912546833829
$ws.Range("E2").Value = 9.869214871283734
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.780652580904117
$ws.Range("K2").Value = 16.88971373102335
$ws.Range("L2").Value = 10.19306339876419
$ws.Range("M2").Value = 17.95538585512775
$ws.Range("N2").Value = 26.31208829589515

$ws.Range("B3").Value = 20.42907201227153
$ws.Range("C3").Value = 3.91594703968447
$ws.Range("E3").Value = 9.884288074501368
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.7837907481538
$ws.Range("K3").Value = 16.81226962451622
$ws.Range("L3").Value = 10.20477477605824
$ws.Range("M3").Value = 17.95306707913999
$ws.Range("N3").Value = 26.33793667670805

$ws.Range("B4").Value = 20.36210992324326
$ws.Range("C4").Value = 3.784298793765827
$ws.Range("E4").Value = 9.894451095813626
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.785817638587847
$ws.Range("K4").Value = 16.76853815194037
$ws.Range("L4").Value = 10.21336237205222
$ws.Range("M4").Value = 17.95480703291594
$ws.Range("N4").Value = 26.35552557178055

$ws.Range("B5").Value = 20.33602939430222
$ws.Range("C5").Value = 3.729066130111713
$ws.Range("E5").Value = 9.898821322807953
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.786668859426888
$ws.Range("K5").Value = 16.75169090944246
$ws.Range("L5").Value = 10.21721331548569
$ws.Range("M5").Value = 17.9563127936848
$ws.Range("N5").Value = 26.36312498087111

$ws.Range("B6").Value = 20.33177225970848
$ws.Range("C6").Value = 3.719799975284016
$ws.Range("E6").Value = 9.899560820647864
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.786811731468743
$ws.Range("K6").Value = 16.74895260631894
$ws.Range("L6").Value = 10.21787399203838
$ws.Range("M6").Value = 17.95661097156062
$ws.Range("N6").Value = 26.36441293059614

$ws.Range("B7").Value = 20.36175327676901
$ws.Range("C7").Value = 3.783560283657588
$ws.Range("E7").Value = 9.894509107653754
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.785829016109354
$ws.Range("K7").Value = 16.76830698526197
$ws.Range("L7").Value = 10.21341288407422
$ws.Range("M7").Value = 17.95482411297882
$ws.Range("N7").Value = 26.35562631207249

$ws.Range("B8").Value = 20.50459066075792
$ws.Range("C8").Value = 4.050931087711968
$ws.Range("E8").Value = 9.874223886053841
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.781713911251083
$ws.Range("K8").Value = 16.8622266054066
$ws.Range("L8").Value = 10.19681171675575
$ws.Range("M8").Value = 17.95393051993284
$ws.Range("N8").Value = 26.32064416008816

$ws.Range("B9").Value = 20.82101308743575
$ws.Range("C9").Value = 4.523406670466844
$ws.Range("E9").Value = 9.841632365036142
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.774433896804005
$ws.Range("K9").Value = 17.07606671768997
$ws.Range("L9").Value = 10.1753297668309
$ws.Range("M9").Value = 17.97721693888031
$ws.Range("N9").Value = 26.26568386297682

$ws.Range("B10").Value = 21.07408574796472
$ws.Range("C10").Value = 4.837950620878673
$ws.Range("E10").Value = 9.822046439261102
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.769560942700664
$ws.Range("K10").Value = 17.25031604916285
$ws.Range("L10").Value = 10.16628264767652
$ws.Range("M10").Value = 18.00948670124379
$ws.Range("N10").Value = 26.23363288367631

$ws.Range("B11").Value = 21.19332722464831
$ws.Range("C11").Value = 4.973831150373456
$ws.Range("E11").Value = 9.814077936865523
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.767446171777591
$ws.Range("K11").Value = 17.33308382382245
$ws.Range("L11").Value = 10.16362532901947
$ws.Range("M11").Value = 18.02742999511438
$ws.Range("N11").Value = 26.22086336324024

$ws.Range("B12").Value = 21.23904070273553
$ws.Range("C12").Value = 5.024240956938788
$ws.Range("E12").Value = 9.811195415178474
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.766659931180399
$ws.Range("K12").Value = 17.36490873305801
$ws.Range("L12").Value = 10.1628283016076
$ws.Range("M12").Value = 18.03469081832289
$ws.Range("N12").Value = 26.21628846932314

$ws.Range("B13").Value = 21.2291711440084
$ws.Range("C13").Value = 5.013430906068853
$ws.Range("E13").Value = 9.811810220610687
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.766828615174184
$ws.Range("K13").Value = 17.35803354958351
$ws.Range("L13").Value = 10.16299065753143
$ws.Range("M13").Value = 18.03310639590314
$ws.Range("N13").Value = 26.21726215786497

$ws.Range("B14").Value = 21.19707708174818
$ws.Range("C14").Value = 4.977999385637158
$ws.Range("E14").Value = 9.813838087068131
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.767381195623032
$ws.Range("K14").Value = 17.33569252993762
$ws.Range("L14").Value = 10.16355556640201
$ws.Range("M14").Value = 18.02801802606039
$ws.Range("N14").Value = 26.22048175864741

$ws.Range("B15").Value = 21.17749039162855
$ws.Range("C15").Value = 4.956160216025872
$ws.Range("E15").Value = 9.815097780757176
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.767721562992534
$ws.Range("K15").Value = 17.32207021790115
$ws.Range("L15").Value = 10.16392882429321
$ws.Range("M15").Value = 18.02496184917377
$ws.Range("N15").Value = 26.22248780912247

$ws.Range("B16").Value = 21.06637311401879
$ws.Range("C16").Value = 4.828924677221901
$ws.Range("E16").Value = 9.822586108050277
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.769701192420932
$ws.Range("K16").Value = 17.24497566336962
$ws.Range("L16").Value = 10.16648561614811
$ws.Range("M16").Value = 18.00837947052877
$ws.Range("N16").Value = 26.23450386045345

$ws.Range("B17").Value = 20.99923787274719
$ws.Range("C17").Value = 4.749016874373043
$ws.Range("E17").Value = 9.827420764270698
$ws.Range("F17").Value = 20.20408069617459
$ws.Range("G17").Value = 3.770941685157272
$ws.Range("K17").Value = 17.19856274941753
$ws.Range("L17").Value = 10.16842731633626
$ws.Range("M17").Value = 17.99904043714672
$ws.Range("N17").Value = 26.24233926580072

$ws.Range("B18").Value = 20.96101304663819
$ws.Range("C18").Value = 4.70237802440548
$ws.Range("E18").Value = 9.830290146502611
$ws.Range("F18").Value = 19.95656407809808
$ws.Range("G18").Value = 3.771664785784516
$ws.Range("K18").Value = 17.17219809839073
$ws.Range("L18").Value = 10.16968142039918
$ws.Range("M18").Value = 17.99397627528864
$ws.Range("N18").Value = 26.24701639259055

$ws.Range("B19").Value = 20.94813867034717
$ws.Range("C19").Value = 4.68647077123947
$ws.Range("E19").Value = 9.831276901103701
$ws.Range("F19").Value = 19.87204792380562
$ws.Range("G19").Value = 3.77191126696906
$ws.Range("K19").Value = 17.16332892193008
$ws.Range("L19").Value = 10.17012963094391
$ws.Range("M19").Value = 17.9923145258307
$ws.Range("N19").Value = 26.24862924657545

$ws.Range("B20").Value = 21.00634444690668
$ws.Range("C20").Value = 4.757593437302785
$ws.Range("E20").Value = 9.826896938131416
$ws.Range("F20").Value = 20.24955283636157
$ws.Range("G20").Value = 3.770808639432305
$ws.Range("K20").Value = 17.20346939217711
$ws.Range("L20").Value = 10.1682064132198
$ws.Range("M20").Value = 18.00000280033593
$ws.Range("N20").Value = 26.2414875349216

$ws.Range("B21").Value = 21.20648897226614
$ws.Range("C21").Value = 4.988434907975273
$ws.Range("E21").Value = 9.813238793270418
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.767218494398906
$ws.Range("K21").Value = 17.34224170341858
$ws.Range("L21").Value = 10.16338396414134
$ws.Range("M21").Value = 18.02949998161143
$ws.Range("N21").Value = 26.21952900778541

$ws.Range("B22").Value = 21.34053924500842
$ws.Range("C22").Value = 5.133208309662484
$ws.Range("E22").Value = 9.805098981159858
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.764957055586394
$ws.Range("K22").Value = 17.43573981846555
$ws.Range("L22").Value = 10.16145158083702
$ws.Range("M22").Value = 18.05149305713548
$ws.Range("N22").Value = 26.20669721246666

$ws.Range("B23").Value = 21.26870832277119
$ws.Range("C23").Value = 5.056500062636947
$ws.Range("E23").Value = 9.809371504180717
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.766156284883245
$ws.Range("K23").Value = 17.38558881379089
$ws.Range("L23").Value = 10.16237152392258
$ws.Range("M23").Value = 18.03950767710933
$ws.Range("N23").Value = 26.21340666985983

$ws.Range("B24").Value = 21.0031304031448
$ws.Range("C24").Value = 4.753718149129639
$ws.Range("E24").Value = 9.827133479965132
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.770868758483073
$ws.Range("K24").Value = 17.20125010453411
$ws.Range("L24").Value = 10.16830585427646
$ws.Range("M24").Value = 17.99956676574413
$ws.Range("N24").Value = 26.24187206527746

$ws.Range("B25").Value = 20.73167732949463
$ws.Range("C25").Value = 4.401254679566022
$ws.Range("E25").Value = 9.84968207583465
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.776319387187605
$ws.Range("K25").Value = 17.01513418340499
$ws.Range("L25").Value = 10.17995706918104
$ws.Range("M25").Value = 17.96824597802708
$ws.Range("N25").Value = 26.27909047383874
